$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025947876515066
$ws.Range("D2").Value = 1.030600807718336
$ws.Range("E2").Value = 1.050410729326888
$ws.Range("F2").Value = 1.055248167865698
$ws.Range("I2").Value = 1.034427613586252
$ws.Range("J2").Value = 1.031114052051779
$ws.Range("K2").Value = 1.033411348109769
$ws.Range("L2").Value = 1.053165082743529
$ws.Range("M2").Value = 1.05798914817972
$ws.Range("N2").Value = 1.014418142644217

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026666382971807
$ws.Range("D3").Value = 1.0311165512366
$ws.Range("E3").Value = 1.051413850396665
$ws.Range("F3").Value = 1.056253265434431
$ws.Range("I3").Value = 1.034565360641614
$ws.Range("J3").Value = 1.031473300402063
$ws.Range("K3").Value = 1.033736266708319
$ws.Range("L3").Value = 1.053980153303403
$ws.Range("M3").Value = 1.058807160298519
$ws.Range("N3").Value = 1.014537384150392

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027131849328898
$ws.Range("D4").Value = 1.031450720881166
$ws.Range("E4").Value = 1.052064285558571
$ws.Range("F4").Value = 1.056904836685662
$ws.Range("I4").Value = 1.034653585377006
$ws.Range("J4").Value = 1.03170560278864
$ws.Range("K4").Value = 1.033946244942257
$ws.Range("L4").Value = 1.054508292801576
$ws.Range("M4").Value = 1.059337059869827
$ws.Range("N4").Value = 1.014614475555416

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027327659939537
$ws.Range("D5").Value = 1.031591311554151
$ws.Range("E5").Value = 1.052338050052227
$ws.Range("F5").Value = 1.057179044498202
$ws.Range("I5").Value = 1.034690457301479
$ws.Range("J5").Value = 1.03180322438865
$ws.Range("K5").Value = 1.034034455003585
$ws.Range("L5").Value = 1.05473049726121
$ws.Range("M5").Value = 1.059559969848356
$ws.Range("N5").Value = 1.014646868621032

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027360544887877
$ws.Range("D6").Value = 1.031614923496963
$ws.Range("E6").Value = 1.052384035132376
$ws.Range("F6").Value = 1.057225101999744
$ws.Range("I6").Value = 1.034696635471475
$ws.Range("J6").Value = 1.031819613204398
$ws.Range("K6").Value = 1.034049262028386
$ws.Range("L6").Value = 1.05476781657277
$ws.Range("M6").Value = 1.059597405598378
$ws.Range("N6").Value = 1.014652306599869

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027134465257029
$ws.Range("D7").Value = 1.031452599047299
$ws.Range("E7").Value = 1.052067942351329
$ws.Range("F7").Value = 1.05690849953694
$ws.Range("I7").Value = 1.034654078917912
$ws.Range("J7").Value = 1.031706907365252
$ws.Range("K7").Value = 1.033947423865431
$ws.Range("L7").Value = 1.054511261223807
$ws.Range("M7").Value = 1.059340037854601
$ws.Range("N7").Value = 1.014614908457153

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026190585264151
$ws.Range("D8").Value = 1.030775011607819
$ws.Range("E8").Value = 1.05074945871169
$ws.Range("F8").Value = 1.055587594734494
$ws.Range("I8").Value = 1.034474353214678
$ws.Range("J8").Value = 1.031235493188343
$ws.Range("K8").Value = 1.033521210221983
$ws.Range("L8").Value = 1.053440386888278
$ws.Range("M8").Value = 1.05826547604773
$ws.Range("N8").Value = 1.014458454256484

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024531600444338
$ws.Range("D9").Value = 1.029584533330137
$ws.Range("E9").Value = 1.048436511458466
$ws.Range("F9").Value = 1.053269291429451
$ws.Range("I9").Value = 1.034150737628164
$ws.Range("J9").Value = 1.030403662016392
$ws.Range("K9").Value = 1.032768182201057
$ws.Range("L9").Value = 1.051559046402411
$ws.Range("M9").Value = 1.056376544885564
$ws.Range("N9").Value = 1.014182275643378

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023428579052497
$ws.Range("D10").Value = 1.028793351200279
$ws.Range("E10").Value = 1.046901611621445
$ws.Range("F10").Value = 1.051730094127848
$ws.Range("I10").Value = 1.033930383499079
$ws.Range("J10").Value = 1.029848410034663
$ws.Range("K10").Value = 1.032264897442128
$ws.Range("L10").Value = 1.050308706915235
$ws.Range("M10").Value = 1.05512041657066
$ws.Range("N10").Value = 1.013997853630154

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022951685838917
$ws.Range("D11").Value = 1.028451368762979
$ws.Range("E11").Value = 1.046238675423942
$ws.Range("F11").Value = 1.05106512595046
$ws.Range("I11").Value = 1.033833882537338
$ws.Range("J11").Value = 1.029607828652474
$ws.Range("K11").Value = 1.032046683800485
$ws.Range("L11").Value = 1.049768232938604
$ws.Range("M11").Value = 1.054577264574291
$ws.Range("N11").Value = 1.013917930090273

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022774656859983
$ws.Range("D12").Value = 1.02832443391635
$ws.Range("E12").Value = 1.04599268596476
$ws.Range("F12").Value = 1.050818356110725
$ws.Range("I12").Value = 1.033797875258975
$ws.Range("J12").Value = 1.029518444096269
$ws.Range("K12").Value = 1.031965587436845
$ws.Range("L12").Value = 1.04956761792305
$ws.Range("M12").Value = 1.054375629179576
$ws.Range("N12").Value = 1.013888233183918

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022812625151264
$ws.Range("D13").Value = 1.028351657649796
$ws.Range("E13").Value = 1.046045439997979
$ws.Range("F13").Value = 1.05087127869564
$ws.Range("I13").Value = 1.033805606286374
$ws.Range("J13").Value = 1.029537618360276
$ws.Range("K13").Value = 1.031982984768522
$ws.Range("L13").Value = 1.049610644121683
$ws.Range("M13").Value = 1.054418875411649
$ws.Range("N13").Value = 1.013894603705551

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022937050305454
$ws.Range("D14").Value = 1.028440874388293
$ws.Range("E14").Value = 1.046218336657837
$ws.Range("F14").Value = 1.051044723199766
$ws.Range("I14").Value = 1.033830909476102
$ws.Range("J14").Value = 1.029600440537614
$ws.Range("K14").Value = 1.032039981201855
$ws.Range("L14").Value = 1.049751647153973
$ws.Range("M14").Value = 1.05456059496728
$ws.Range("N14").Value = 1.013915475531159

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023013727433827
$ws.Range("D15").Value = 1.028495856097637
$ws.Range("E15").Value = 1.046324897767956
$ws.Range("F15").Value = 1.051151618463887
$ws.Range("I15").Value = 1.033846478090109
$ws.Range("J15").Value = 1.029639144484647
$ws.Range("K15").Value = 1.032075093041726
$ws.Range("L15").Value = 1.049838542505977
$ws.Range("M15").Value = 1.054647928398264
$ws.Range("N15").Value = 1.013928334072705

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023460244262937
$ws.Range("D16").Value = 1.028816060342687
$ws.Range("E16").Value = 1.04694564413935
$ws.Range("F16").Value = 1.051774257904227
$ws.Range("I16").Value = 1.03393676512207
$ws.Range("J16").Value = 1.029864373491245
$ws.Range("K16").Value = 1.032279373594557
$ws.Range("L16").Value = 1.050344596113697
$ws.Range("M16").Value = 1.055156479898293
$ws.Range("N16").Value = 1.014003156503197

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023740527182634
$ws.Range("D17").Value = 1.029017079170469
$ws.Range("E17").Value = 1.047335474144239
$ws.Range("F17").Value = 1.052165229630685
$ws.Range("I17").Value = 1.03399310945905
$ws.Range("J17").Value = 1.030005613389438
$ws.Range("K17").Value = 1.032407437104003
$ws.Range("L17").Value = 1.050662280324783
$ws.Range("M17").Value = 1.055475685155424
$ws.Range("N17").Value = 1.014050072830441

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023904081036266
$ws.Range("D18").Value = 1.029134388269791
$ws.Range("E18").Value = 1.047563017933602
$ws.Range("F18").Value = 1.052393422846753
$ws.Range("I18").Value = 1.034025869313559
$ws.Range("J18").Value = 1.030087981333639
$ws.Range("K18").Value = 1.032482106514238
$ws.Range("L18").Value = 1.050847669925735
$ws.Range("M18").Value = 1.055661945409939
$ws.Range("N18").Value = 1.014077431781349

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023959860416812
$ws.Range("D19").Value = 1.029174397460584
$ws.Range("E19").Value = 1.047640632015394
$ws.Range("F19").Value = 1.052471255561103
$ws.Range("I19").Value = 1.034037021773931
$ws.Range("J19").Value = 1.030116064111654
$ws.Range("K19").Value = 1.032507562090941
$ws.Range("L19").Value = 1.050910898180702
$ws.Range("M19").Value = 1.055725467725279
$ws.Range("N19").Value = 1.014086759353289

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023710448268613
$ws.Range("D20").Value = 1.028995505717804
$ws.Range("E20").Value = 1.047293632251876
$ws.Range("F20").Value = 1.052123266950438
$ws.Range("I20").Value = 1.033987075084676
$ws.Range("J20").Value = 1.029990461217737
$ws.Range("K20").Value = 1.032393699975512
$ws.Range("L20").Value = 1.050628186521115
$ws.Range("M20").Value = 1.055441429868635
$ws.Range("N20").Value = 1.014045039826979

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022900407136299
$ws.Range("D21").Value = 1.028414599710118
$ws.Range("E21").Value = 1.046167415863016
$ws.Range("F21").Value = 1.050993641784692
$ws.Range("I21").Value = 1.033823462798555
$ws.Range("J21").Value = 1.029581941567956
$ws.Range("K21").Value = 1.032023198325642
$ws.Range("L21").Value = 1.049710121371331
$ws.Range("M21").Value = 1.05451885889337
$ws.Range("N21").Value = 1.013909329564457

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022391741760162
$ws.Range("D22").Value = 1.028049898029101
$ws.Range("E22").Value = 1.045460792397
$ws.Range("F22").Value = 1.050284726947573
$ws.Range("I22").Value = 1.033719653271877
$ws.Range("J22").Value = 1.029324962806759
$ws.Range("K22").Value = 1.031790006227184
$ws.Range("L22").Value = 1.049133714226052
$ws.Range("M22").Value = 1.053939470446578
$ws.Range("N22").Value = 1.013823946956634

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022661333645047
$ws.Range("D23").Value = 1.02824318178581
$ws.Range("E23").Value = 1.045835246753177
$ws.Range("F23").Value = 1.050660410011139
$ws.Range("I23").Value = 1.033774773545346
$ws.Range("J23").Value = 1.029461203689837
$ws.Range("K23").Value = 1.03191364841439
$ws.Range("L23").Value = 1.049439200739061
$ws.Range("M23").Value = 1.054246551425188
$ws.Range("N23").Value = 1.013869215081031

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023724039420417
$ws.Range("D24").Value = 1.029005253652806
$ws.Range("E24").Value = 1.047312538299873
$ws.Range("F24").Value = 1.052142227628339
$ws.Range("I24").Value = 1.033989802082744
$ws.Range("J24").Value = 1.029997307878097
$ws.Range("K24").Value = 1.032399907278831
$ws.Range("L24").Value = 1.050643591766167
$ws.Range("M24").Value = 1.0554569081321
$ws.Range("N24").Value = 1.014047314044889

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024959972585346
$ws.Range("D25").Value = 1.029891872820123
$ws.Range("E25").Value = 1.049033224374845
$ws.Range("F25").Value = 1.053867517539299
$ws.Range("I25").Value = 1.034235215525556
$ws.Range("J25").Value = 1.030618838296126
$ws.Range("K25").Value = 1.032963086168237
$ws.Range("L25").Value = 1.052044738010013
$ws.Range("M25").Value = 1.056864327715711
$ws.Range("N25").Value = 1.014253729587672

Write-Output "updated vm_pu values for 380 kV case"